$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing cell values per diff ---
$ws.Cells.Item(60, 17).Value = 0    # Q60: 1 -> 0
$ws.Cells.Item(69, 17).Value = 0    # Q69: 1 -> 0
$ws.Cells.Item(76, 17).Value = 0    # Q76: 2 -> 0
$ws.Cells.Item(370, 15).Value = 2   # O370: 0 -> 2
$ws.Cells.Item(372, 18).Value = 0   # R372: (blank) -> 0
$ws.Cells.Item(373, 18).Value = 0   # R373: (blank) -> 0

# --- Append new weekly rows 374:396 ---
$dateFormat = $ws.Cells.Item(373, 1).NumberFormat

# Row 374
$ws.Cells.Item(374, 1).Value = 45474
$ws.Cells.Item(374, 2).Value = 281.3500061035156
$ws.Cells.Item(374, 3).Value = 339.4500122070312
$ws.Cells.Item(374, 4).Value = 277
$ws.Cells.Item(374, 5).Value = 328.4500122070312
$ws.Cells.Item(374, 6).Value = 325.0246276855469
$ws.Cells.Item(374, 7).Value = 231877798
$ws.Cells.Item(374, 8).Value = 2024
$ws.Cells.Item(374, 9).Value = 7
$ws.Cells.Item(374, 10).Value = 1
$ws.Cells.Item(374, 11).Value = 0
$ws.Cells.Item(374, 12).Value = 0
$ws.Cells.Item(374, 13).Value = 0
$ws.Cells.Item(374, 14).Value = 27
$ws.Cells.Item(374, 15).Value = 0
$ws.Cells.Item(374, 16).Value = 0
$ws.Cells.Item(374, 17).Value = 0
$ws.Cells.Item(374, 1).NumberFormat = $dateFormat

# Row 375
$ws.Cells.Item(375, 1).Value = 45481
$ws.Cells.Item(375, 2).Value = 330.0499877929688
$ws.Cells.Item(375, 3).Value = 353.7000122070312
$ws.Cells.Item(375, 4).Value = 319
$ws.Cells.Item(375, 5).Value = 334
$ws.Cells.Item(375, 6).Value = 330.5167236328125
$ws.Cells.Item(375, 7).Value = 160053196
$ws.Cells.Item(375, 8).Value = 2024
$ws.Cells.Item(375, 9).Value = 7
$ws.Cells.Item(375, 10).Value = 8
$ws.Cells.Item(375, 11).Value = 0
$ws.Cells.Item(375, 12).Value = 0
$ws.Cells.Item(375, 13).Value = 0
$ws.Cells.Item(375, 14).Value = 28
$ws.Cells.Item(375, 15).Value = 1
$ws.Cells.Item(375, 16).Value = 0
$ws.Cells.Item(375, 17).Value = 0
$ws.Cells.Item(375, 1).NumberFormat = $dateFormat

# Row 376
$ws.Cells.Item(376, 1).Value = 45488
$ws.Cells.Item(376, 2).Value = 336.4500122070312
$ws.Cells.Item(376, 3).Value = 346.7999877929688
$ws.Cells.Item(376, 4).Value = 309
$ws.Cells.Item(376, 5).Value = 318
$ws.Cells.Item(376, 6).Value = 314.68359375
$ws.Cells.Item(376, 7).Value = 93000590
$ws.Cells.Item(376, 8).Value = 2024
$ws.Cells.Item(376, 9).Value = 7
$ws.Cells.Item(376, 10).Value = 15
$ws.Cells.Item(376, 11).Value = 0
$ws.Cells.Item(376, 12).Value = 0
$ws.Cells.Item(376, 13).Value = 0
$ws.Cells.Item(376, 14).Value = 29
$ws.Cells.Item(376, 15).Value = 0
$ws.Cells.Item(376, 16).Value = 0
$ws.Cells.Item(376, 17).Value = 1
$ws.Cells.Item(376, 1).NumberFormat = $dateFormat

# Row 377
$ws.Cells.Item(377, 1).Value = 45495
$ws.Cells.Item(377, 2).Value = 317
$ws.Cells.Item(377, 3).Value = 329.3999938964844
$ws.Cells.Item(377, 4).Value = 271.1499938964844
$ws.Cells.Item(377, 5).Value = 314.25
$ws.Cells.Item(377, 6).Value = 310.9726867675781
$ws.Cells.Item(377, 7).Value = 196040984
$ws.Cells.Item(377, 8).Value = 2024
$ws.Cells.Item(377, 9).Value = 7
$ws.Cells.Item(377, 10).Value = 22
$ws.Cells.Item(377, 11).Value = 0
$ws.Cells.Item(377, 12).Value = 0
$ws.Cells.Item(377, 13).Value = 0
$ws.Cells.Item(377, 14).Value = 30
$ws.Cells.Item(377, 15).Value = 2
$ws.Cells.Item(377, 16).Value = 0
$ws.Cells.Item(377, 17).Value = 0
$ws.Cells.Item(377, 1).NumberFormat = $dateFormat

# Row 378
$ws.Cells.Item(378, 1).Value = 45502
$ws.Cells.Item(378, 2).Value = 316.0499877929688
$ws.Cells.Item(378, 3).Value = 321.3500061035156
$ws.Cells.Item(378, 4).Value = 296.5
$ws.Cells.Item(378, 5).Value = 306.2000122070312
$ws.Cells.Item(378, 6).Value = 303.0066528320312
$ws.Cells.Item(378, 7).Value = 48968555
$ws.Cells.Item(378, 8).Value = 2024
$ws.Cells.Item(378, 9).Value = 7
$ws.Cells.Item(378, 10).Value = 29
$ws.Cells.Item(378, 11).Value = 0
$ws.Cells.Item(378, 12).Value = 0
$ws.Cells.Item(378, 13).Value = 0
$ws.Cells.Item(378, 14).Value = 31
$ws.Cells.Item(378, 15).Value = 0
$ws.Cells.Item(378, 16).Value = 0
$ws.Cells.Item(378, 17).Value = 0
$ws.Cells.Item(378, 1).NumberFormat = $dateFormat

# Row 379
$ws.Cells.Item(379, 1).Value = 45509
$ws.Cells.Item(379, 2).Value = 291.1499938964844
$ws.Cells.Item(379, 3).Value = 304.1499938964844
$ws.Cells.Item(379, 4).Value = 284.1000061035156
$ws.Cells.Item(379, 5).Value = 292.2000122070312
$ws.Cells.Item(379, 6).Value = 289.1526794433594
$ws.Cells.Item(379, 7).Value = 60136899
$ws.Cells.Item(379, 8).Value = 2024
$ws.Cells.Item(379, 9).Value = 8
$ws.Cells.Item(379, 10).Value = 5
$ws.Cells.Item(379, 11).Value = 0
$ws.Cells.Item(379, 12).Value = 0
$ws.Cells.Item(379, 13).Value = 0
$ws.Cells.Item(379, 14).Value = 32
$ws.Cells.Item(379, 15).Value = 0
$ws.Cells.Item(379, 16).Value = 0
$ws.Cells.Item(379, 17).Value = 0
$ws.Cells.Item(379, 1).NumberFormat = $dateFormat

# Row 380
$ws.Cells.Item(380, 1).Value = 45516
$ws.Cells.Item(380, 2).Value = 295
$ws.Cells.Item(380, 3).Value = 312.8999938964844
$ws.Cells.Item(380, 4).Value = 285
$ws.Cells.Item(380, 5).Value = 291.7999877929688
$ws.Cells.Item(380, 6).Value = 288.7568054199219
$ws.Cells.Item(380, 7).Value = 59053408
$ws.Cells.Item(380, 8).Value = 2024
$ws.Cells.Item(380, 9).Value = 8
$ws.Cells.Item(380, 10).Value = 12
$ws.Cells.Item(380, 11).Value = 0
$ws.Cells.Item(380, 12).Value = 0
$ws.Cells.Item(380, 13).Value = 0
$ws.Cells.Item(380, 14).Value = 33
$ws.Cells.Item(380, 15).Value = 0
$ws.Cells.Item(380, 16).Value = 0
$ws.Cells.Item(380, 17).Value = 0
$ws.Cells.Item(380, 1).NumberFormat = $dateFormat

# Row 381
$ws.Cells.Item(381, 1).Value = 45523
$ws.Cells.Item(381, 2).Value = 292
$ws.Cells.Item(381, 3).Value = 297.7999877929688
$ws.Cells.Item(381, 4).Value = 282.1499938964844
$ws.Cells.Item(381, 5).Value = 284.3500061035156
$ws.Cells.Item(381, 6).Value = 281.384521484375
$ws.Cells.Item(381, 7).Value = 31209476
$ws.Cells.Item(381, 8).Value = 2024
$ws.Cells.Item(381, 9).Value = 8
$ws.Cells.Item(381, 10).Value = 19
$ws.Cells.Item(381, 11).Value = 0
$ws.Cells.Item(381, 12).Value = 0
$ws.Cells.Item(381, 13).Value = 0
$ws.Cells.Item(381, 14).Value = 34
$ws.Cells.Item(381, 15).Value = 0
$ws.Cells.Item(381, 16).Value = 0
$ws.Cells.Item(381, 17).Value = 2
$ws.Cells.Item(381, 1).NumberFormat = $dateFormat

# Row 382
$ws.Cells.Item(382, 1).Value = 45530
$ws.Cells.Item(382, 2).Value = 286
$ws.Cells.Item(382, 3).Value = 299.7000122070312
$ws.Cells.Item(382, 4).Value = 277
$ws.Cells.Item(382, 5).Value = 280.6499938964844
$ws.Cells.Item(382, 6).Value = 277.7231140136719
$ws.Cells.Item(382, 7).Value = 87269518
$ws.Cells.Item(382, 8).Value = 2024
$ws.Cells.Item(382, 9).Value = 8
$ws.Cells.Item(382, 10).Value = 26
$ws.Cells.Item(382, 11).Value = 0
$ws.Cells.Item(382, 12).Value = 0
$ws.Cells.Item(382, 13).Value = 0
$ws.Cells.Item(382, 14).Value = 35
$ws.Cells.Item(382, 15).Value = 0
$ws.Cells.Item(382, 16).Value = 0
$ws.Cells.Item(382, 17).Value = 0
$ws.Cells.Item(382, 1).NumberFormat = $dateFormat

# Row 383
$ws.Cells.Item(383, 1).Value = 45537
$ws.Cells.Item(383, 2).Value = 279.9500122070312
$ws.Cells.Item(383, 3).Value = 281.7999877929688
$ws.Cells.Item(383, 4).Value = 253.1000061035156
$ws.Cells.Item(383, 5).Value = 254.1000061035156
$ws.Cells.Item(383, 6).Value = 251.4500122070312
$ws.Cells.Item(383, 7).Value = 37571502
$ws.Cells.Item(383, 8).Value = 2024
$ws.Cells.Item(383, 9).Value = 9
$ws.Cells.Item(383, 10).Value = 2
$ws.Cells.Item(383, 11).Value = 0
$ws.Cells.Item(383, 12).Value = 0
$ws.Cells.Item(383, 13).Value = 0
$ws.Cells.Item(383, 14).Value = 36
$ws.Cells.Item(383, 15).Value = 0
$ws.Cells.Item(383, 16).Value = 0
$ws.Cells.Item(383, 17).Value = 1
$ws.Cells.Item(383, 1).NumberFormat = $dateFormat

# Row 384
$ws.Cells.Item(384, 1).Value = 45544
$ws.Cells.Item(384, 2).Value = 253.9499969482422
$ws.Cells.Item(384, 3).Value = 258.6000061035156
$ws.Cells.Item(384, 4).Value = 245
$ws.Cells.Item(384, 5).Value = 252.4499969482422
$ws.Cells.Item(384, 6).Value = 249.8171997070312
$ws.Cells.Item(384, 7).Value = 43748013
$ws.Cells.Item(384, 8).Value = 2024
$ws.Cells.Item(384, 9).Value = 9
$ws.Cells.Item(384, 10).Value = 9
$ws.Cells.Item(384, 11).Value = 0
$ws.Cells.Item(384, 12).Value = 0
$ws.Cells.Item(384, 13).Value = 0
$ws.Cells.Item(384, 14).Value = 37
$ws.Cells.Item(384, 15).Value = 0
$ws.Cells.Item(384, 16).Value = 0
$ws.Cells.Item(384, 17).Value = 0
$ws.Cells.Item(384, 1).NumberFormat = $dateFormat

# Row 385
$ws.Cells.Item(385, 1).Value = 45551
$ws.Cells.Item(385, 2).Value = 253.5
$ws.Cells.Item(385, 3).Value = 254.6999969482422
$ws.Cells.Item(385, 4).Value = 228.0500030517578
$ws.Cells.Item(385, 5).Value = 251.0500030517578
$ws.Cells.Item(385, 6).Value = 251.0500030517578
$ws.Cells.Item(385, 7).Value = 68791769
$ws.Cells.Item(385, 8).Value = 2024
$ws.Cells.Item(385, 9).Value = 9
$ws.Cells.Item(385, 10).Value = 16
$ws.Cells.Item(385, 11).Value = 0
$ws.Cells.Item(385, 12).Value = 0
$ws.Cells.Item(385, 13).Value = 0
$ws.Cells.Item(385, 14).Value = 38
$ws.Cells.Item(385, 15).Value = 0
$ws.Cells.Item(385, 16).Value = 0
$ws.Cells.Item(385, 17).Value = 0
$ws.Cells.Item(385, 1).NumberFormat = $dateFormat

# Row 386
$ws.Cells.Item(386, 1).Value = 45558
$ws.Cells.Item(386, 2).Value = 251.25
$ws.Cells.Item(386, 3).Value = 252.8500061035156
$ws.Cells.Item(386, 4).Value = 232.1499938964844
$ws.Cells.Item(386, 5).Value = 238.4499969482422
$ws.Cells.Item(386, 6).Value = 238.4499969482422
$ws.Cells.Item(386, 7).Value = 34012591
$ws.Cells.Item(386, 8).Value = 2024
$ws.Cells.Item(386, 9).Value = 9
$ws.Cells.Item(386, 10).Value = 23
$ws.Cells.Item(386, 11).Value = 0
$ws.Cells.Item(386, 12).Value = 0
$ws.Cells.Item(386, 13).Value = 0
$ws.Cells.Item(386, 14).Value = 39
$ws.Cells.Item(386, 15).Value = 0
$ws.Cells.Item(386, 16).Value = 0
$ws.Cells.Item(386, 17).Value = 0
$ws.Cells.Item(386, 1).NumberFormat = $dateFormat

# Row 387
$ws.Cells.Item(387, 1).Value = 45565
$ws.Cells.Item(387, 2).Value = 238
$ws.Cells.Item(387, 3).Value = 243.1000061035156
$ws.Cells.Item(387, 4).Value = 220
$ws.Cells.Item(387, 5).Value = 224.3300018310547
$ws.Cells.Item(387, 6).Value = 224.3300018310547
$ws.Cells.Item(387, 7).Value = 20172014
$ws.Cells.Item(387, 8).Value = 2024
$ws.Cells.Item(387, 9).Value = 9
$ws.Cells.Item(387, 10).Value = 30
$ws.Cells.Item(387, 11).Value = 0
$ws.Cells.Item(387, 12).Value = 0
$ws.Cells.Item(387, 13).Value = 0
$ws.Cells.Item(387, 14).Value = 40
$ws.Cells.Item(387, 15).Value = 0
$ws.Cells.Item(387, 16).Value = 0
$ws.Cells.Item(387, 17).Value = 0
$ws.Cells.Item(387, 1).NumberFormat = $dateFormat

# Row 388
$ws.Cells.Item(388, 1).Value = 45572
$ws.Cells.Item(388, 2).Value = 224.3200073242188
$ws.Cells.Item(388, 3).Value = 232.5599975585938
$ws.Cells.Item(388, 4).Value = 206.5599975585938
$ws.Cells.Item(388, 5).Value = 223.1000061035156
$ws.Cells.Item(388, 6).Value = 223.1000061035156
$ws.Cells.Item(388, 7).Value = 31546912
$ws.Cells.Item(388, 8).Value = 2024
$ws.Cells.Item(388, 9).Value = 10
$ws.Cells.Item(388, 10).Value = 7
$ws.Cells.Item(388, 11).Value = 0
$ws.Cells.Item(388, 12).Value = 0
$ws.Cells.Item(388, 13).Value = 0
$ws.Cells.Item(388, 14).Value = 41
$ws.Cells.Item(388, 15).Value = 0
$ws.Cells.Item(388, 16).Value = 0
$ws.Cells.Item(388, 17).Value = 0
$ws.Cells.Item(388, 1).NumberFormat = $dateFormat

# Row 389
$ws.Cells.Item(389, 1).Value = 45579
$ws.Cells.Item(389, 2).Value = 223.1000061035156
$ws.Cells.Item(389, 3).Value = 228.5200042724609
$ws.Cells.Item(389, 4).Value = 209.8999938964844
$ws.Cells.Item(389, 5).Value = 215.0599975585938
$ws.Cells.Item(389, 6).Value = 215.0599975585938
$ws.Cells.Item(389, 7).Value = 13047830
$ws.Cells.Item(389, 8).Value = 2024
$ws.Cells.Item(389, 9).Value = 10
$ws.Cells.Item(389, 10).Value = 14
$ws.Cells.Item(389, 11).Value = 0
$ws.Cells.Item(389, 12).Value = 0
$ws.Cells.Item(389, 13).Value = 0
$ws.Cells.Item(389, 14).Value = 42
$ws.Cells.Item(389, 15).Value = 0
$ws.Cells.Item(389, 16).Value = 0
$ws.Cells.Item(389, 17).Value = 0
$ws.Cells.Item(389, 1).NumberFormat = $dateFormat

# Row 390
$ws.Cells.Item(390, 1).Value = 45586
$ws.Cells.Item(390, 2).Value = 216.1499938964844
$ws.Cells.Item(390, 3).Value = 217.5500030517578
$ws.Cells.Item(390, 4).Value = 193
$ws.Cells.Item(390, 5).Value = 193.9900054931641
$ws.Cells.Item(390, 6).Value = 193.9900054931641
$ws.Cells.Item(390, 7).Value = 26378353
$ws.Cells.Item(390, 8).Value = 2024
$ws.Cells.Item(390, 9).Value = 10
$ws.Cells.Item(390, 10).Value = 21
$ws.Cells.Item(390, 11).Value = 0
$ws.Cells.Item(390, 12).Value = 0
$ws.Cells.Item(390, 13).Value = 0
$ws.Cells.Item(390, 14).Value = 43
$ws.Cells.Item(390, 15).Value = 0
$ws.Cells.Item(390, 16).Value = 0
$ws.Cells.Item(390, 17).Value = 0
$ws.Cells.Item(390, 1).NumberFormat = $dateFormat

# Row 391
$ws.Cells.Item(391, 1).Value = 45593
$ws.Cells.Item(391, 2).Value = 194.9900054931641
$ws.Cells.Item(391, 3).Value = 223.8999938964844
$ws.Cells.Item(391, 4).Value = 192.1100006103516
$ws.Cells.Item(391, 5).Value = 219.3399963378906
$ws.Cells.Item(391, 6).Value = 219.3399963378906
$ws.Cells.Item(391, 7).Value = 33235914
$ws.Cells.Item(391, 8).Value = 2024
$ws.Cells.Item(391, 9).Value = 10
$ws.Cells.Item(391, 10).Value = 28
$ws.Cells.Item(391, 11).Value = 0
$ws.Cells.Item(391, 12).Value = 0
$ws.Cells.Item(391, 13).Value = 0
$ws.Cells.Item(391, 14).Value = 44
$ws.Cells.Item(391, 15).Value = 2
$ws.Cells.Item(391, 16).Value = 0
$ws.Cells.Item(391, 17).Value = 0
$ws.Cells.Item(391, 1).NumberFormat = $dateFormat

# Row 392
$ws.Cells.Item(392, 1).Value = 45600
$ws.Cells.Item(392, 2).Value = 219.6999969482422
$ws.Cells.Item(392, 3).Value = 231.7400054931641
$ws.Cells.Item(392, 4).Value = 213.7299957275391
$ws.Cells.Item(392, 5).Value = 218.5200042724609
$ws.Cells.Item(392, 6).Value = 218.5200042724609
$ws.Cells.Item(392, 7).Value = 21698188
$ws.Cells.Item(392, 8).Value = 2024
$ws.Cells.Item(392, 9).Value = 11
$ws.Cells.Item(392, 10).Value = 4
$ws.Cells.Item(392, 11).Value = 0
$ws.Cells.Item(392, 12).Value = 0
$ws.Cells.Item(392, 13).Value = 0
$ws.Cells.Item(392, 14).Value = 45
$ws.Cells.Item(392, 15).Value = 0
$ws.Cells.Item(392, 16).Value = 0
$ws.Cells.Item(392, 17).Value = 0
$ws.Cells.Item(392, 1).NumberFormat = $dateFormat

# Row 393
$ws.Cells.Item(393, 1).Value = 45607
$ws.Cells.Item(393, 2).Value = 215.9799957275391
$ws.Cells.Item(393, 3).Value = 224.0299987792969
$ws.Cells.Item(393, 4).Value = 200
$ws.Cells.Item(393, 5).Value = 202.0399932861328
$ws.Cells.Item(393, 6).Value = 202.0399932861328
$ws.Cells.Item(393, 7).Value = 12434038
$ws.Cells.Item(393, 8).Value = 2024
$ws.Cells.Item(393, 9).Value = 11
$ws.Cells.Item(393, 10).Value = 11
$ws.Cells.Item(393, 11).Value = 0
$ws.Cells.Item(393, 12).Value = 0
$ws.Cells.Item(393, 13).Value = 0
$ws.Cells.Item(393, 14).Value = 46
$ws.Cells.Item(393, 15).Value = 0
$ws.Cells.Item(393, 16).Value = 0
$ws.Cells.Item(393, 17).Value = 0
$ws.Cells.Item(393, 1).NumberFormat = $dateFormat

# Row 394
$ws.Cells.Item(394, 1).Value = 45614
$ws.Cells.Item(394, 2).Value = 203
$ws.Cells.Item(394, 3).Value = 214.4700012207031
$ws.Cells.Item(394, 4).Value = 198.0099945068359
$ws.Cells.Item(394, 5).Value = 205.9799957275391
$ws.Cells.Item(394, 6).Value = 205.9799957275391
$ws.Cells.Item(394, 7).Value = 12412113
$ws.Cells.Item(394, 8).Value = 2024
$ws.Cells.Item(394, 9).Value = 11
$ws.Cells.Item(394, 10).Value = 18
$ws.Cells.Item(394, 11).Value = 0
$ws.Cells.Item(394, 12).Value = 0
$ws.Cells.Item(394, 13).Value = 0
$ws.Cells.Item(394, 14).Value = 47
$ws.Cells.Item(394, 15).Value = 0
$ws.Cells.Item(394, 16).Value = 0
$ws.Cells.Item(394, 17).Value = 0
$ws.Cells.Item(394, 1).NumberFormat = $dateFormat

# Row 395
$ws.Cells.Item(395, 1).Value = 45621
$ws.Cells.Item(395, 2).Value = 215
$ws.Cells.Item(395, 3).Value = 242.5
$ws.Cells.Item(395, 4).Value = 210.8500061035156
$ws.Cells.Item(395, 5).Value = 238.6100006103516
$ws.Cells.Item(395, 6).Value = 238.6100006103516
$ws.Cells.Item(395, 7).Value = 65373016
$ws.Cells.Item(395, 8).Value = 2024
$ws.Cells.Item(395, 9).Value = 11
$ws.Cells.Item(395, 10).Value = 25
$ws.Cells.Item(395, 11).Value = 0
$ws.Cells.Item(395, 12).Value = 0
$ws.Cells.Item(395, 13).Value = 0
$ws.Cells.Item(395, 14).Value = 48
$ws.Cells.Item(395, 15).Value = 0
$ws.Cells.Item(395, 16).Value = 0
$ws.Cells.Item(395, 17).Value = 0
$ws.Cells.Item(395, 1).NumberFormat = $dateFormat

# Row 396
$ws.Cells.Item(396, 1).Value = 45628
$ws.Cells.Item(396, 2).Value = 240
$ws.Cells.Item(396, 3).Value = 249.8999938964844
$ws.Cells.Item(396, 4).Value = 234.1999969482422
$ws.Cells.Item(396, 5).Value = 247.7700042724609
$ws.Cells.Item(396, 6).Value = 247.7700042724609
$ws.Cells.Item(396, 7).Value = 47928937
$ws.Cells.Item(396, 8).Value = 2024
$ws.Cells.Item(396, 9).Value = 12
$ws.Cells.Item(396, 10).Value = 2
$ws.Cells.Item(396, 11).Value = 0
$ws.Cells.Item(396, 12).Value = 0
$ws.Cells.Item(396, 13).Value = 0
$ws.Cells.Item(396, 14).Value = 49
$ws.Cells.Item(396, 15).Value = 0
$ws.Cells.Item(396, 16).Value = 0
$ws.Cells.Item(396, 17).Value = 0
$ws.Cells.Item(396, 1).NumberFormat = $dateFormat

